# UPDATED THE NEW TEMPLATE FORMAT (MAJOR CHANGE)
# Rewrites the tag syntax used in the "Overlay" sheet's Content/Params/PreProcess
# columns from the old shorthand (<!T>/<!F>/<Process=.../<PreProcess=...)
# to the new explicit Key=Value syntax (<Type=...>/<File=...>/<Function=...),
# introduces a dedicated "AddSpace" function row for the name/TIN concat
# separators, and widens column C to fit the longer strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overlay")

# --- Column C (Content) text-syntax updates ---------------------------------
$ws.Range("C2").Value = "<Type=Text><Text=2024-10-25>"
$ws.Range("C3").Value = "<Type=Text><Text=2023/2024>"
$ws.Range("C4").Value = "<Type=Text><Text=5249087539>"
$ws.Range("C5").Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=B>"
$ws.Range("C6").Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=D>"

# Row 7: the literal-space placeholder becomes an empty Text tag, and the
# "add a space" behaviour moves into a new Function tag in column E.
$ws.Range("C7").Value = "<Type=Text><Text=>"
$ws.Range("E7").Value = "<Function=AddSpace(None)>"
$ws.Range("E7").NumberFormat = "@"

$ws.Range("C8").Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=C>"

# Row 9: same AddSpace pattern as row 7 (separator before the TIN literal).
$ws.Range("C9").Value = "<Type=Text><Text=>"
$ws.Range("E9").Value = "<Function=AddSpace(None)>"
$ws.Range("E9").NumberFormat = "@"

$ws.Range("C10").Value = "<Type=Text><Text=TIN>"
$ws.Range("C11").Value = "<Type=File><File=EMP01.xlsx><Sheet=PERSONAL DATA><PrimeryKey=A><Value=C>"
$ws.Range("C12").Value = "<Type=File><File=EMP01.xlsx><Sheet=PERSONAL DATA><PrimeryKey=A><Value=B>"

$ws.Range("C13").Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=E>"
$ws.Range("D13").Value = "<X=170><Y=200><Font=Helvetica><FontSize=12><LineSpace=1.2X><Function=SrinkToFit(300,2)>"
$ws.Range("E13").Value = "<Function=NumberToText(text)>"

$ws.Range("C14").Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=E>"
$ws.Range("E14").Value = "<Function=NumberToCurrency(text,USD,2)>"

# --- Widen column C to fit the longer tag strings ---------------------------
$ws.Columns.Item(3).ColumnWidth = 62
